# Auto-generated cell updates for cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.699.82"
$ws.Range("D3").Value = "1.894.66"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  -1.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4885"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3787"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07326"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9118"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07658"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.881.72"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.478"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.623"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008775"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("D20").Value = "28.060.58"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.121"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "2.101.00"
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("E25").Value = "  -2.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.162"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.881"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.200"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.32%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7677"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("E35").Value = "  -0.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02035"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.544"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.97%  "
$ws.Range("E38").Value = "  -2.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05275"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5466"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.886"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.540"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "112.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.05%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1519"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4780"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9999"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.638"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06052"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.19%  "
